# Remove the last two rows (old rows 11 and 12) and the old row 10, then
# rebuild row 10 with the new "proceedings-article:None" / "includes HTML" data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = "MathML"
$ws.Range("C2").Value = 6

# Update row 5
$ws.Range("B5").Value = "U+2009 thin space from general punc; "
$ws.Range("C5").Value = 2

# Update row 6
$ws.Range("B6").Value = "includes HTML"

# Update row 7
$ws.Range("B7").Value = "includes disclosures; includes figure"

# Update row 8
$ws.Range("B8").Value = "includes hyperlinks"

# Update row 9
$ws.Range("B9").Value = "includes tex-math, not MathML"

# Remove old rows 10, 11, 12 (originally "includes script from the landing
# page", "non-Latin char; MathML", "proceedings-article:None" / "less than
# and greater than markup")
$ws.Rows(12).Delete()
$ws.Rows(11).Delete()
$ws.Rows(10).Delete()

# Rebuild row 10 with the new data, copying formatting from row 9 so the
# cell styles match what a real edit in Excel would look like.
$ws.Cells.Item(9,1).Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$ws.Cells.Item(9,2).Copy()
$ws.Cells.Item(10,2).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(10,1).Value = "proceedings-article:None"
$ws.Cells.Item(10,2).Value = "includes HTML"
$ws.Cells.Item(10,3).Value = 1
